# Apply the "Add files via upload" edit to Fragen.xlsx (Tabelle1 sheet)
#
# 1) Row 11, column C: question text changed from the generic
#    "A3.4 Vermietung land-/forstw. Betriebsmittel (an wen? Bezirk? Kondition – upload)"
#    to the new, more specific
#    "A3.4 Vermietung land-/forstw. Betriebsmittel (an Nichtlandwirte?)"
#
# 2) Column D ("Attribute") on rows 8,9,10,11,12,34,35,36,37,38: the generic
#    attribute string "Beschreibung:text;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
#    is replaced with a new attribute string that asks for revenue info:
#    "Bitte geben Sie hier Ihre Umsätze bekannt:info;Umsatz:text;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
#
# 3) The sheet's selected cell moved to D39 (view scroll position is not
#    something this host persists across save/load, only the active
#    selection is).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

$newQuestion = "A3.4 Vermietung land-/forstw. Betriebsmittel (an Nichtlandwirte?)"
$newAttribute = "Bitte geben Sie hier Ihre Umsätze bekannt:info;Umsatz:text;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"

# Row 11 question text (column C)
$ws.Range("C11").Value = $newQuestion

# Column D ("Attribute") on the rows that referenced the generic
# "Beschreibung:text;..." string
$attributeRows = @(8, 9, 10, 11, 12, 34, 35, 36, 37, 38)
foreach ($r in $attributeRows) {
    $ws.Cells.Item($r, 4).Value = $newAttribute
}

# Update selection to match the new view state
$ws.Range("D39").Select() | Out-Null
